# Update "want to go" counts (column F) on the 展览 and 全部类型 sheets
$wb = $excel.ActiveWorkbook

$updates = @(
    @{Row=4;  Value=4},
    @{Row=7;  Value=1780},
    @{Row=11; Value=2089},
    @{Row=13; Value=145},
    @{Row=14; Value=1344},
    @{Row=16; Value=23},
    @{Row=25; Value=1129},
    @{Row=27; Value=337},
    @{Row=29; Value=273},
    @{Row=30; Value=332}
)

$ws1 = $wb.Worksheets.Item("展览")
foreach ($u in $updates) {
    $ws1.Cells.Item($u.Row, 6).Value = $u.Value
}

$updates4 = @(
    @{Row=4;  Value=4},
    @{Row=7;  Value=1780},
    @{Row=12; Value=2089},
    @{Row=14; Value=145},
    @{Row=15; Value=1344},
    @{Row=17; Value=23},
    @{Row=26; Value=1129},
    @{Row=28; Value=337},
    @{Row=30; Value=273},
    @{Row=31; Value=332}
)

$ws4 = $wb.Worksheets.Item("全部类型")
foreach ($u in $updates4) {
    $ws4.Cells.Item($u.Row, 6).Value = $u.Value
}
